# Auto-generated: apply 2024-08-30 daily crime-count increments
# to the 2024 (column K) figures across the Citywide Totals,
# By Neighborhood summary, and every per-neighborhood sheet.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{}

$sheetUpdates['Citywide Totals'] = @{ 'K2' = 5282; 'K3' = 5448; 'K4' = 1130; 'K5' = 388; 'K6' = 6064; 'K7' = 18312 }
$sheetUpdates['By Neighborhood'] = @{ 'K2' = 161; 'K7' = 541; 'K8' = 1232; 'K10' = 101; 'K11' = 350; 'K14' = 96; 'K15' = 187; 'K19' = 538; 'K20' = 423; 'K23' = 193; 'K27' = 175; 'K29' = 980; 'K30' = 71; 'K33' = 785; 'K34' = 104; 'K37' = 617; 'K38' = 17; 'K42' = 677; 'K43' = 163; 'K47' = 124; 'K48' = 230; 'K52' = 477; 'K54' = 360; 'K60' = 114; 'K63' = 52; 'K64' = 117; 'K65' = 416; 'K66' = 60; 'K67' = 694; 'K73' = 157; 'K76' = 253; 'K78' = 210; 'K79' = 450; 'K83' = 406; 'K85' = 864; 'K89' = 265; 'K90' = 165; 'K91' = 201; 'K94' = 241; 'K97' = 146; 'K98' = 86; 'K99' = 309; 'K101' = 18312 }
$sheetUpdates['Bridgeport'] = @{ 'K6' = 36; 'K7' = 96 }
$sheetUpdates['Auburn Gresham'] = @{ 'K2' = 184; 'K6' = 144; 'K7' = 541 }
$sheetUpdates['Belmont Cragin'] = @{ 'K4' = 21; 'K7' = 350 }
$sheetUpdates['Uptown'] = @{ 'K4' = 31; 'K7' = 265 }
$sheetUpdates['South Shore'] = @{ 'K2' = 289; 'K3' = 292; 'K7' = 864 }
$sheetUpdates['Little Village'] = @{ 'K2' = 127; 'K3' = 134; 'K7' = 477 }
$sheetUpdates['Austin'] = @{ 'K2' = 338; 'K3' = 369; 'K7' = 1232 }
$sheetUpdates['South Chicago'] = @{ 'K2' = 136; 'K7' = 406 }
$sheetUpdates['Garfield Park'] = @{ 'K2' = 215; 'K3' = 290; 'K6' = 228; 'K7' = 785 }
$sheetUpdates['Grand Crossing'] = @{ 'K2' = 177; 'K6' = 180; 'K7' = 617 }
$sheetUpdates['New City'] = @{ 'K2' = 133; 'K6' = 154; 'K7' = 416 }
$sheetUpdates['Woodlawn'] = @{ 'K3' = 128; 'K7' = 309 }
$sheetUpdates['Fuller Park'] = @{ 'K2' = 20; 'K7' = 71 }
$sheetUpdates['North Lawndale'] = @{ 'K2' = 199; 'K3' = 247; 'K7' = 694 }
$sheetUpdates['Loop'] = @{ 'K3' = 90; 'K6' = 192; 'K7' = 360 }
$sheetUpdates['Englewood'] = @{ 'K2' = 280; 'K3' = 354; 'K6' = 271; 'K7' = 980 }
$sheetUpdates['Lake View'] = @{ 'K6' = 114; 'K7' = 230 }
$sheetUpdates['Chatham'] = @{ 'K2' = 158; 'K3' = 169; 'K7' = 538 }
$sheetUpdates['River North'] = @{ 'K2' = 55; 'K7' = 253 }
$sheetUpdates['Humboldt Park'] = @{ 'K2' = 179; 'K6' = 256; 'K7' = 677 }
$sheetUpdates['Avondale'] = @{ 'K3' = 17; 'K6' = 50; 'K7' = 101 }
$sheetUpdates['Rogers Park'] = @{ 'K3' = 47; 'K7' = 210 }
$sheetUpdates['Douglas'] = @{ 'K2' = 53; 'K6' = 53; 'K7' = 193 }
$sheetUpdates['Washington Park'] = @{ 'K2' = 50; 'K3' = 96; 'K7' = 201 }
$sheetUpdates['Roseland'] = @{ 'K6' = 110; 'K7' = 450 }
$sheetUpdates['Near South Side'] = @{ 'K6' = 45; 'K7' = 117 }
$sheetUpdates['Chicago Lawn'] = @{ 'K2' = 139; 'K5' = 8; 'K7' = 423 }
$sheetUpdates['Garfield Ridge'] = @{ 'K3' = 28; 'K7' = 104 }
$sheetUpdates['West Loop'] = @{ 'K6' = 104; 'K7' = 241 }
$sheetUpdates['Kenwood'] = @{ 'K4' = 11; 'K7' = 124 }
$sheetUpdates['Brighton Park'] = @{ 'K2' = 65; 'K3' = 46; 'K7' = 187 }
$sheetUpdates['Wicker Park'] = @{ 'K6' = 53; 'K7' = 86 }
$sheetUpdates['North Center'] = @{ 'K3' = 13; 'K7' = 60 }
$sheetUpdates['Portage Park'] = @{ 'K2' = 51; 'K7' = 157 }
$sheetUpdates['Albany Park'] = @{ 'K6' = 51; 'K7' = 161 }
$sheetUpdates['West Town'] = @{ 'K6' = 86; 'K7' = 146 }
$sheetUpdates['Edgewater'] = @{ 'K2' = 48; 'K3' = 42; 'K7' = 175 }
$sheetUpdates['Washington Heights'] = @{ 'K2' = 59; 'K7' = 165 }
$sheetUpdates['Morgan Park'] = @{ 'K6' = 33; 'K7' = 114 }
$sheetUpdates['Hyde Park'] = @{ 'K4' = 23; 'K6' = 65; 'K7' = 163 }
$sheetUpdates['Grant Park'] = @{ 'K5' = 8; 'K6' = 17 }

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellValues = $sheetUpdates[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range($cellRef).Value = $cellValues[$cellRef]
    }
}

Write-Host "Updated $($sheetUpdates.Keys.Count) sheets."
